$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet used to list the email in column A and the name in column B.
# Swap the row-2 values: column A now holds the (hyperlinked) e-mail
# address, column B holds the plain name.
$ws.Range("B2").Value = "Saw Myint Win"
$ws.Range("A2").Value = "sawmyintwin@gmail.com"

# Turn the e-mail address in A2 into a clickable mailto: link. Excel
# auto-applies the built-in "Hyperlink" cell style (underline + theme
# colour) the first time a hyperlink is added to the workbook.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:sawmyintwin@gmail.com")

# Move the active selection, matching the saved cursor position.
$ws.Range("A20").Select()

$wb.Save()
